# FE and BE integration
#
# The template used split runs to bracket replaceable tokens, e.g.
#   <change>RELEASE_NO | <change>
#   <change> | TITLE | <change> | " "
#   <change> | RELEASE_NO | <change>
# split across multiple w:r runs (so FE/BE code scanning for the literal
# "<change>...<change>" token across run text wouldn't see it as one
# string). Re-typing the same text via Find/Replace merges each matched
# range into a single run (using the formatting of the first run in the
# match), collapsing the split runs without touching anything else.

$d = $word.ActiveDocument

# 1) "<change>RELEASE_NO" + "<change>"  ->  one run "<change>RELEASE_NO<change>"
$d.Content.Find.Execute("<change>RELEASE_NO<change>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<change>RELEASE_NO<change>", 2) | Out-Null

# 2) "<change>" + "TITLE" + "<change>" + " "  ->  one run "<change>TITLE<change> "
$d.Content.Find.Execute("<change>TITLE<change> ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<change>TITLE<change> ", 2) | Out-Null

# 3) "<change>" + "RELEASE_NO" + "<change>"  ->  one run "<change>RELEASE_NO<change>"
#    (the "Parte fixa unde mai variaza numarul release-ului: " occurrence)
$d.Content.Find.Execute("<change>RELEASE_NO<change>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<change>RELEASE_NO<change>", 2) | Out-Null
